# Apply updated symbol list values (price & 1h volume %) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected cells to remain plain text so numeric-looking
# strings (prices, percentages) are preserved exactly as authored,
# instead of being reinterpreted as numbers by Excel.
$cellRefs = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D20", "E20", "D21", "E21", "D22", "E22", "E23", "D24", "E24", "D25", "E25", "E26", "D27", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51"
)
foreach ($ref in $cellRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "332.23"
$ws.Range("E2").Value = "0.06%"
$ws.Range("D3").Value = "41.16"
$ws.Range("E3").Value = "0.01%"
$ws.Range("D4").Value = "5.695"
$ws.Range("E4").Value = "-0.58%"
$ws.Range("D5").Value = "0.08417"
$ws.Range("E5").Value = "3.10%"
$ws.Range("D6").Value = "8.826"
$ws.Range("E6").Value = "0.96%"
$ws.Range("D7").Value = "4.513"
$ws.Range("E7").Value = "0.10%"
$ws.Range("D8").Value = "1.979"
$ws.Range("E8").Value = "-3.88%"
$ws.Range("D10").Value = "0.9270"
$ws.Range("E10").Value = "0.66%"
$ws.Range("D11").Value = "0.1247"
$ws.Range("E11").Value = "0.51%"
$ws.Range("D12").Value = "0.1972"
$ws.Range("E12").Value = "1.54%"
$ws.Range("D13").Value = "0.09362"
$ws.Range("E13").Value = "-0.78%"
$ws.Range("D14").Value = "0.03968"
$ws.Range("E14").Value = "8.61%"
$ws.Range("D15").Value = "0.1063"
$ws.Range("E15").Value = "0.83%"
$ws.Range("D16").Value = "0.001306"
$ws.Range("E16").Value = "0.43%"
$ws.Range("D17").Value = "0.006102"
$ws.Range("E17").Value = "-1.39%"
$ws.Range("D18").Value = "3.433"
$ws.Range("E18").Value = "1.45%"
$ws.Range("D20").Value = "9.062"
$ws.Range("E20").Value = "8.94%"
$ws.Range("D21").Value = "0.1374"
$ws.Range("E21").Value = "-2.96%"
$ws.Range("D22").Value = "0.2632"
$ws.Range("E22").Value = "-0.62%"
$ws.Range("E23").Value = "0.17%"
$ws.Range("D24").Value = "0.001244"
$ws.Range("E24").Value = "-1.17%"
$ws.Range("D25").Value = "0.004372"
$ws.Range("E25").Value = "0.68%"
$ws.Range("E26").Value = "-3.89%"
$ws.Range("D27").Value = "0.0003999"
$ws.Range("E27").Value = "0.15%"
$ws.Range("D39").Value = "0.02811"
$ws.Range("E39").Value = "1.04%"
$ws.Range("D40").Value = "0.05533"
$ws.Range("E40").Value = "0.44%"
$ws.Range("D41").Value = "0.007918"
$ws.Range("E41").Value = "4.08%"
$ws.Range("E42").Value = "0.95%"
$ws.Range("D43").Value = "0.008973"
$ws.Range("E43").Value = "-9.72%"
$ws.Range("D44").Value = "0.002095"
$ws.Range("E44").Value = "-1.10%"
$ws.Range("D45").Value = "0.01012"
$ws.Range("E45").Value = "-14.89%"
$ws.Range("D46").Value = "0.00007175"
$ws.Range("E46").Value = "6.19%"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").Value = "0.31%"
$ws.Range("D48").Value = "0.003460"
$ws.Range("E48").Value = "18.42%"
$ws.Range("D49").Value = "0.002283"
$ws.Range("E49").Value = "0.25%"
$ws.Range("D50").Value = "0.00002105"
$ws.Range("E50").Value = "0.31%"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").Value = "0.31%"
